# Auto-generated script applying scheduled market-data update to Sheets workbook.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H:N) per leve row
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1109.579
$ws.Range("I19").Value = 1169.5
$ws.Range("J19").Value = 1006.8571
$ws.Range("K19").Value = 1169.5
$ws.Range("L19").Value = 1006.8571
$ws.Range("M19").Value = -994.5
$ws.Range("N19").Value = -1356.8571
$ws.Range("H113").Value = 2821.625
$ws.Range("I113").Value = 2251
$ws.Range("K113").Value = 2251
$ws.Range("M113").Value = 1003
$ws.Range("H132").Value = 4204.6733
$ws.Range("I132").Value = 4251.9556
$ws.Range("K132").Value = 12755.8668
$ws.Range("M132").Value = -10225.8668
$ws.Range("H133").Value = 110000
$ws.Range("J133").Value = 110000
$ws.Range("L133").Value = 110000
$ws.Range("N133").Value = -120120
$ws.Range("H135").Value = 1498.1765
$ws.Range("I135").Value = 1339.1
$ws.Range("K135").Value = 12051.9
$ws.Range("M135").Value = -9516.9
$ws.Range("H137").Value = 48677.055
$ws.Range("I137").Value = 51334.53
$ws.Range("J137").Value = 3500
$ws.Range("K137").Value = 154003.59
$ws.Range("L137").Value = 10500
$ws.Range("M137").Value = -151453.59
$ws.Range("N137").Value = -15600

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15081578
$ws.Range("I32").Value = 14920364
$ws.Range("K32").Value = 14920364
$ws.Range("M32").Value = -14920077
$ws.Range("H45").Value = 6483.846
$ws.Range("I45").Value = 6429
$ws.Range("K45").Value = 6429
$ws.Range("M45").Value = -6052
$ws.Range("H55").Value = 43997.5
$ws.Range("H61").Value = 2584.647
$ws.Range("I61").Value = 2691.2173
$ws.Range("J61").Value = 2361.818
$ws.Range("K61").Value = 2691.2173
$ws.Range("L61").Value = 2361.818
$ws.Range("M61").Value = -2479.2173
$ws.Range("N61").Value = -2785.818
$ws.Range("H74").Value = 2323.6667
$ws.Range("I74").Value = 2625.44
$ws.Range("J74").Value = 1784.7858
$ws.Range("K74").Value = 2625.44
$ws.Range("L74").Value = 1784.7858
$ws.Range("M74").Value = -1751.44
$ws.Range("N74").Value = -3532.7858
$ws.Range("H77").Value = 2323.6667
$ws.Range("I77").Value = 2625.44
$ws.Range("J77").Value = 1784.7858
$ws.Range("K77").Value = 13127.2
$ws.Range("L77").Value = 8923.929
$ws.Range("M77").Value = -8759.200000000001
$ws.Range("N77").Value = -17659.929
$ws.Range("H132").Value = 3485.8076
$ws.Range("I132").Value = 2759.5
$ws.Range("K132").Value = 8278.5
$ws.Range("M132").Value = -5748.5
$ws.Range("H136").Value = 2584.647
$ws.Range("I136").Value = 2691.2173
$ws.Range("J136").Value = 2361.818
$ws.Range("K136").Value = 8073.651899999999
$ws.Range("L136").Value = 7085.454000000001
$ws.Range("M136").Value = -5523.651899999999
$ws.Range("N136").Value = -12185.454

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4205732.5
$ws.Range("I134").Value = 5955769
$ws.Range("J134").Value = 5645.8
$ws.Range("K134").Value = 17867307
$ws.Range("L134").Value = 16937.4
$ws.Range("M134").Value = -17864772
$ws.Range("N134").Value = -22007.4

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1364
$ws.Range("J16").Value = 2237.6667
$ws.Range("L16").Value = 2237.6667
$ws.Range("N16").Value = -2811.6667
$ws.Range("H22").Value = 434.84616
$ws.Range("I22").Value = 355.2857
$ws.Range("K22").Value = 355.2857
$ws.Range("M22").Value = -5.28570000000002
$ws.Range("H31").Value = 6919.4707
$ws.Range("I31").Value = 2214.2
$ws.Range("K31").Value = 2214.2
$ws.Range("M31").Value = -1919.2
$ws.Range("H34").Value = 6919.4707
$ws.Range("I34").Value = 2214.2
$ws.Range("K34").Value = 2214.2
$ws.Range("M34").Value = -2012.2
$ws.Range("H76").Value = 7978.6
$ws.Range("I76").Value = 7978.6
$ws.Range("K76").Value = 7978.6
$ws.Range("M76").Value = -7663.6
$ws.Range("H79").Value = 7978.6
$ws.Range("I79").Value = 7978.6
$ws.Range("K79").Value = 7978.6
$ws.Range("M79").Value = -6886.6
$ws.Range("H99").Value = 4079.2
$ws.Range("I99").Value = 3799
$ws.Range("K99").Value = 3799
$ws.Range("M99").Value = -2301
$ws.Range("H113").Value = 1364
$ws.Range("J113").Value = 2237.6667
$ws.Range("L113").Value = 2237.6667
$ws.Range("N113").Value = -6577.6667
$ws.Range("H126").Value = 4079.2
$ws.Range("I126").Value = 3799
$ws.Range("K126").Value = 11397
$ws.Range("M126").Value = -8927
$ws.Range("H132").Value = 18029.508
$ws.Range("I132").Value = 20106.104
$ws.Range("K132").Value = 60318.312
$ws.Range("M132").Value = -57788.312
$ws.Range("H134").Value = 44906.082
$ws.Range("I134").Value = 69912
$ws.Range("J134").Value = 3229.5557
$ws.Range("K134").Value = 209736
$ws.Range("L134").Value = 9688.667099999999
$ws.Range("M134").Value = -207201
$ws.Range("N134").Value = -14758.6671
$ws.Range("H138").Value = 94543.05499999999
$ws.Range("J138").Value = 97280.88
$ws.Range("L138").Value = 97280.88
$ws.Range("N138").Value = -107560.88
$ws.Range("H58").Value = 3303.1875
$ws.Range("I58").Value = 3303.1875
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 3303.1875
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -3100.1875
$ws.Range("N58").ClearContents()
$ws.Range("H136").Value = 3303.1875
$ws.Range("I136").Value = 3303.1875
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 9909.5625
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -7359.5625
$ws.Range("N136").ClearContents()

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 70560660
$ws.Range("I4").Value = 60283012
$ws.Range("J4").Value = 88832030
$ws.Range("K4").Value = 180849036
$ws.Range("L4").Value = 266496090
$ws.Range("M4").Value = -180848924
$ws.Range("N4").Value = -266496314
$ws.Range("H68").Value = 965.3333
$ws.Range("I68").Value = 898
$ws.Range("J68").Value = 1100
$ws.Range("K68").Value = 2694
$ws.Range("L68").Value = 3300
$ws.Range("M68").Value = -1883
$ws.Range("N68").Value = -4922
$ws.Range("H71").Value = 965.3333
$ws.Range("I71").Value = 898
$ws.Range("J71").Value = 1100
$ws.Range("K71").Value = 8082
$ws.Range("L71").Value = 9900
$ws.Range("M71").Value = -4026
$ws.Range("N71").Value = -18012
$ws.Range("H129").Value = 1781.6875
$ws.Range("I129").Value = 858.2222
$ws.Range("J129").Value = 2969
$ws.Range("K129").Value = 2574.6666
$ws.Range("L129").Value = 8907
$ws.Range("M129").Value = 2425.3334
$ws.Range("N129").Value = -18907
$ws.Range("H130").Value = 3567
$ws.Range("I130").Value = 2492.25
$ws.Range("K130").Value = 7476.75
$ws.Range("M130").Value = -2456.75
$ws.Range("H131").Value = 1471.2858
$ws.Range("I131").Value = 717.2222
$ws.Range("J131").Value = 2036.8334
$ws.Range("K131").Value = 2151.6666
$ws.Range("L131").Value = 6110.5002
$ws.Range("M131").Value = 2888.3334
$ws.Range("N131").Value = -16190.5002

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1006.6316
$ws.Range("I97").Value = 847
$ws.Range("J97").Value = 1713.5714
$ws.Range("K97").Value = 847
$ws.Range("L97").Value = 1713.5714
$ws.Range("M97").Value = -351
$ws.Range("N97").Value = -2705.5714
$ws.Range("H122").Value = 2075.6667
$ws.Range("I122").Value = 1939.8
$ws.Range("K122").Value = 5819.4
$ws.Range("M122").Value = -3369.4
$ws.Range("H126").Value = 3499.6
$ws.Range("I126").Value = 3249
$ws.Range("J126").Value = 3666.6667
$ws.Range("K126").Value = 9747
$ws.Range("L126").Value = 11000.0001
$ws.Range("M126").Value = -7277
$ws.Range("N126").Value = -15940.0001
$ws.Range("H130").Value = 99999
$ws.Range("J130").Value = 99999
$ws.Range("L130").Value = 99999
$ws.Range("N130").Value = -110039
$ws.Range("H132").Value = 7548.857
$ws.Range("I132").Value = 7772.1
$ws.Range("J132").Value = 6990.75
$ws.Range("K132").Value = 23316.3
$ws.Range("L132").Value = 20972.25
$ws.Range("M132").Value = -20786.3
$ws.Range("N132").Value = -26032.25
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6640.9443
$ws.Range("I7").Value = 6136.75
$ws.Range("K7").Value = 6136.75
$ws.Range("M7").Value = -6024.75
$ws.Range("H40").Value = 37046010
$ws.Range("I40").Value = 55561524
$ws.Range("K40").Value = 55561524
$ws.Range("M40").Value = -55561388
$ws.Range("H122").Value = 22672.842
$ws.Range("I122").Value = 21266.334
$ws.Range("K122").Value = 63799.00199999999
$ws.Range("M122").Value = -61349.00199999999
$ws.Range("H126").Value = 6640.9443
$ws.Range("I126").Value = 6136.75
$ws.Range("K126").Value = 18410.25
$ws.Range("M126").Value = -15940.25
$ws.Range("H132").Value = 30935.703
$ws.Range("I132").Value = 36975.035
$ws.Range("K132").Value = 110925.105
$ws.Range("M132").Value = -108395.105
$ws.Range("H136").Value = 3867.5
$ws.Range("I136").Value = 2747.75
$ws.Range("K136").Value = 8243.25
$ws.Range("M136").Value = -5693.25

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H129").Value = 102831.836
$ws.Range("J129").Value = 102831.836
$ws.Range("L129").Value = 102831.836
$ws.Range("N129").Value = -112831.836
$ws.Range("H136").Value = 34724
$ws.Range("I136").Value = 2474.6667
$ws.Range("J136").Value = 55092
$ws.Range("K136").Value = 7424.000100000001
$ws.Range("L136").Value = 165276
$ws.Range("M136").Value = -4874.000100000001
$ws.Range("N136").Value = -170376

Write-Host "Applied scheduled market data update."